$wb = $excel.ActiveWorkbook

# Sheet "Overview": Latest HO Xliff Generate Date for the 22cf23da row (row 3)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-11-14 06:13:01"

# Sheet "zh-cn": Correspond Handoff Datetime / Correspond Handback DateTime for 22cf23da row (row 3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-11-14 06:12:47"
$wsZhCn.Range("K3").Value = "2016-11-14 06:13:43"

# Sheet "de-de": Correspond Handoff Datetime (same value as Overview's G3) and
# Correspond Handback DateTime for 22cf23da row (row 3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-11-14 06:13:01"
$wsDeDe.Range("K3").Value = "2016-11-14 06:14:01"
